$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 (Trees): Flashcards -> Done, Own Implementation -> Done, Practice Questions -> 2
$ws.Range("E5").Value = "Done"
$ws.Range("F5").Value = "Done"
$ws.Range("G5").Value = 2

# Row 6 (Tries): Notes -> Done, Flashcards -> Done
$ws.Range("D6").Value = "Done"
$ws.Range("E6").Value = "Done"

# Row 7 (Graphs): Flashcards -> Done
$ws.Range("E7").Value = "Done"

# Update selection to match recorded state
[void]$ws.Range("G5:G7").Select()
